$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 6208
    $ws.Range("G6").Value = "不可售"
    $ws.Range("F8").Value = 1871
    $ws.Range("F9").Value = 1401
    $ws.Range("F12").Value = 216
    $ws.Range("F13").Value = 5572
}
